# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column C ("municipio-nombre") was re-classified from an
# iaest-measure/medida/xsd:int triple to an sdmx-dimension/dim/URI-Municipio
# triple, matching the pattern already used by the other "refArea" columns
# (D/provincia-nombre and H/comarca-nombre).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: semantic type for column C -> sdmx-dimension:refArea (was iaest-measure:municipio-nombre)
$ws.Range("C2").Value = "sdmx-dimension:refArea"

# Row 3: dim/medida flag for column C -> dim (was medida)
$ws.Range("C3").Value = "dim"

# Row 4: concept/URI scheme for column C -> URI-Municipio (was xsd:int)
$ws.Range("C4").Value = "URI-Municipio"
